$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: Datum, Zeit, Taetigkeit (copy format from row 21, then update values)
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("A22").Value = 44247
$ws.Range("B22").Value = 4
$ws.Range("D22").Value = "Meeting+Investigating Apexcharts"

# Row 23: Datum, Zeit, Taetigkeit (copy format from row 21, then update values)
$ws.Range("A21").Copy($ws.Range("A23"))
$ws.Range("A23").Value = 44250
$ws.Range("B23").Value = 6
$ws.Range("D23").Value = "Development WebUI"

# Recalculate the shared formulas cascading through C22:C35
$wb.Application.Calculate()

# Update the selected cell shown in the sheet view
$ws.Range("H8").Select() | Out-Null
